$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ccl3"
$ws.Cells.Item(2,3).Value = "Ccr5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 998.7379113333333
$ws.Cells.Item(2,8).Value = 2996.213734
$ws.Cells.Item(2,9).Value = 0.8754681532218018
$ws.Cells.Item(2,10).Value = 0.8754681532218019
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 4.021407666666667
$ws.Cells.Item(2,14).Value = 12.064223
$ws.Cells.Item(2,15).Value = 0.06269882270324605
$ws.Cells.Item(2,16).Value = 0.06269882270324605
$ws.Cells.Item(2,17).Value = 4016.33229362652
$ws.Cells.Item(2,18).Value = 36146.99064263868
$ws.Cells.Item(2,19).Value = 0.054890822521192
$ws.Cells.Item(2,20).Value = 0.05489082252119201

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ccl3"
$ws.Cells.Item(3,3).Value = "Ccr5"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 998.7379113333333
$ws.Cells.Item(3,8).Value = 2996.213734
$ws.Cells.Item(3,9).Value = 0.8754681532218018
$ws.Cells.Item(3,10).Value = 0.8754681532218019
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.4010506666666667
$ws.Cells.Item(3,14).Value = 1.203152
$ws.Cells.Item(3,15).Value = 0.006252886235031953
$ws.Cells.Item(3,16).Value = 0.006252886235031953
$ws.Cells.Item(3,17).Value = 400.5445051655075
$ws.Cells.Item(3,18).Value = 3604.900546489568
$ws.Cells.Item(3,19).Value = 0.005474202764489449
$ws.Cells.Item(3,20).Value = 0.00547420276448945

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ccl3"
$ws.Cells.Item(4,3).Value = "Ccr5"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 998.7379113333333
$ws.Cells.Item(4,8).Value = 2996.213734
$ws.Cells.Item(4,9).Value = 0.8754681532218018
$ws.Cells.Item(4,10).Value = 0.8754681532218019
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 59.71602933333333
$ws.Cells.Item(4,14).Value = 179.148088
$ws.Cells.Item(4,15).Value = 0.931048291061722
$ws.Cells.Item(4,16).Value = 0.931048291061722
$ws.Cells.Item(4,17).Value = 59640.66240949339
$ws.Cells.Item(4,18).Value = 536765.9616854406
$ws.Cells.Item(4,19).Value = 0.8151031279361204
$ws.Cells.Item(4,20).Value = 0.8151031279361205

$ws.Cells.Item(5,1).Value = "M2"
$ws.Cells.Item(5,2).Value = "Ccl3"
$ws.Cells.Item(5,3).Value = "Ccr5"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 142.0664773333333
$ws.Cells.Item(5,8).Value = 426.199432
$ws.Cells.Item(5,9).Value = 0.1245318467781981
$ws.Cells.Item(5,10).Value = 0.1245318467781981
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 4.021407666666667
$ws.Cells.Item(5,14).Value = 12.064223
$ws.Cells.Item(5,15).Value = 0.06269882270324605
$ws.Cells.Item(5,16).Value = 0.06269882270324605
$ws.Cells.Item(5,17).Value = 571.3072211245928
$ws.Cells.Item(5,18).Value = 5141.764990121336
$ws.Cells.Item(5,19).Value = 0.007808000182054047
$ws.Cells.Item(5,20).Value = 0.007808000182054048

$ws.Cells.Item(6,1).Value = "M2"
$ws.Cells.Item(6,2).Value = "Ccl3"
$ws.Cells.Item(6,3).Value = "Ccr5"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 142.0664773333333
$ws.Cells.Item(6,8).Value = 426.199432
$ws.Cells.Item(6,9).Value = 0.1245318467781981
$ws.Cells.Item(6,10).Value = 0.1245318467781981
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.4010506666666667
$ws.Cells.Item(6,14).Value = 1.203152
$ws.Cells.Item(6,15).Value = 0.006252886235031953
$ws.Cells.Item(6,16).Value = 0.006252886235031953
$ws.Cells.Item(6,17).Value = 56.97585544551822
$ws.Cells.Item(6,18).Value = 512.782699009664
$ws.Cells.Item(6,19).Value = 0.0007786834705425034
$ws.Cells.Item(6,20).Value = 0.0007786834705425035

$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Ccl3"
$ws.Cells.Item(7,3).Value = "Ccr5"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 142.0664773333333
$ws.Cells.Item(7,8).Value = 426.199432
$ws.Cells.Item(7,9).Value = 0.1245318467781981
$ws.Cells.Item(7,10).Value = 0.1245318467781981
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 59.71602933333333
$ws.Cells.Item(7,14).Value = 179.148088
$ws.Cells.Item(7,15).Value = 0.931048291061722
$ws.Cells.Item(7,16).Value = 0.931048291061722
$ws.Cells.Item(7,17).Value = 8483.645927720667
$ws.Cells.Item(7,18).Value = 76352.81334948602
$ws.Cells.Item(7,19).Value = 0.1159451631256016
$ws.Cells.Item(7,20).Value = 0.1159451631256016
